$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041052102335306
$ws.Range("D2").Value = 1.046591681889417
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.054335272520966
$ws.Range("I2").Value = 1.040635084502006
$ws.Range("J2").Value = 1.046135393094343
$ws.Range("K2").Value = 1.049356687752546
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.057078766356644
$ws.Range("N2").Value = 1.019267250563567
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041977964249715
$ws.Range("D3").Value = 1.047312737367072
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.055242521938918
$ws.Range("I3").Value = 1.040858210446853
$ws.Range("J3").Value = 1.046707242811987
$ws.Range("K3").Value = 1.049889652212097
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.057798998164449
$ws.Range("N3").Value = 1.019459093332109
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042577365521763
$ws.Range("D4").Value = 1.047779607920522
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.05583031847075
$ws.Range("I4").Value = 1.041001572288545
$ws.Range("J4").Value = 1.047076939447882
$ws.Range("K4").Value = 1.050234143704341
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.058265147480913
$ws.Range("N4").Value = 1.019583063667574
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042829426245993
$ws.Range("D5").Value = 1.047975950649661
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.056077605215432
$ws.Range("I5").Value = 1.041061597955749
$ws.Range("J5").Value = 1.047232280367233
$ws.Range("K5").Value = 1.050378877867628
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.058461142091969
$ws.Range("N5").Value = 1.019635140982401
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042871752569149
$ws.Range("D6").Value = 1.048008921544165
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.05611913609182
$ws.Range("I6").Value = 1.041071662226423
$ws.Range("J6").Value = 1.047258358099364
$ws.Range("K6").Value = 1.050403174061268
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.0584940519071
$ws.Range("N6").Value = 1.019643882647249
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042580733286134
$ws.Range("D7").Value = 1.047782231186936
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.05583362203426
$ws.Range("I7").Value = 1.041002375312479
$ws.Range("J7").Value = 1.047079015434579
$ws.Range("K7").Value = 1.050236078004398
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.058267766272015
$ws.Range("N7").Value = 1.019583759684108
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041364937599648
$ws.Range("D8").Value = 1.046835303157489
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.054641726514636
$ws.Range("I8").Value = 1.040710700924329
$ws.Range("J8").Value = 1.046328719813227
$ws.Range("K8").Value = 1.049536882062334
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.057322148340337
$ws.Range("N8").Value = 1.019332118618856
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03922494317496
$ws.Range("D9").Value = 1.04516904743304
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.052547225248244
$ws.Range("I9").Value = 1.040188978480161
$ws.Range("J9").Value = 1.045004129342662
$ws.Range("K9").Value = 1.048302001873621
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.055656752164114
$ws.Range("N9").Value = 1.018887449295967
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0377999469867
$ws.Range("D10").Value = 1.044059870945031
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.051154852979872
$ws.Range("I10").Value = 1.039835980336373
$ws.Range("J10").Value = 1.044119460627808
$ws.Range("K10").Value = 1.04747691757151
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.054547164756392
$ws.Range("N10").Value = 1.018590187299996
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037183315010028
$ws.Range("D11").Value = 1.043579995629245
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.050552897125374
$ws.Range("I11").Value = 1.039681905215455
$ws.Range("J11").Value = 1.043736019119099
$ws.Range("K11").Value = 1.047119224478105
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.054066876027594
$ws.Range("N11").Value = 1.018461280575282
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036954331542735
$ws.Range("D12").Value = 1.043401810784212
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.050329448065804
$ws.Range("I12").Value = 1.03962449135534
$ws.Range("J12").Value = 1.043593536536785
$ws.Range("K12").Value = 1.046986298228947
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.053888502065185
$ws.Range("N12").Value = 1.01841337063464
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037003446494301
$ws.Range("D13").Value = 1.043440029202114
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.050377372103106
$ws.Range("I13").Value = 1.03963681511924
$ws.Range("J13").Value = 1.043624102024265
$ws.Range("K13").Value = 1.04701481420168
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.053926762661305
$ws.Range("N13").Value = 1.018423648752006
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037164385913501
$ws.Range("D14").Value = 1.043565265537693
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.050534423808091
$ws.Range("I14").Value = 1.039677163108419
$ws.Range("J14").Value = 1.043724242588119
$ws.Range("K14").Value = 1.047108238034967
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.054052131029928
$ws.Range("N14").Value = 1.018457320898779
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037263554125032
$ws.Range("D15").Value = 1.043642436060259
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.050631207687887
$ws.Range("I15").Value = 1.039701998535114
$ws.Range("J15").Value = 1.043785935188533
$ws.Range("K15").Value = 1.047165791207016
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.054129378185396
$ws.Range("N15").Value = 1.018478063687649
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037840879340126
$ws.Range("D16").Value = 1.044091727342388
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.051194822948383
$ws.Range("I16").Value = 1.03984618003496
$ws.Range("J16").Value = 1.044144900570228
$ws.Range("K16").Value = 1.047500647583818
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.054579043613961
$ws.Range("N16").Value = 1.018598738434828
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038203128122681
$ws.Range("D17").Value = 1.044373665357263
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.051548619340032
$ws.Range("I17").Value = 1.039936293761016
$ws.Range("J17").Value = 1.044369970524971
$ws.Range("K17").Value = 1.047710580696985
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.054861153148925
$ws.Range("N17").Value = 1.018674383808737
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038414460317846
$ws.Range("D18").Value = 1.044538154061409
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.051755074411516
$ws.Range("I18").Value = 1.039988737384737
$ws.Range("J18").Value = 1.044501213866149
$ws.Range("K18").Value = 1.04783298995177
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.055025719189122
$ws.Range("N18").Value = 1.018718488089069
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038486525656762
$ws.Range("D19").Value = 1.04459424703926
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.051825485755505
$ws.Range("I19").Value = 1.040006599243789
$ws.Range("J19").Value = 1.044545958294689
$ws.Range("K19").Value = 1.047874721331906
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.055081834690246
$ws.Range("N19").Value = 1.018733523379601
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0381642582566
$ws.Range("D20").Value = 1.044343412030005
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.051510650848901
$ws.Range("I20").Value = 1.039926637632299
$ws.Range("J20").Value = 1.044345826382272
$ws.Range("K20").Value = 1.047688061113366
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.054830883769342
$ws.Range("N20").Value = 1.01866626967528
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037116991571891
$ws.Range("D21").Value = 1.043528384830434
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.050488172005451
$ws.Range("I21").Value = 1.039665286692875
$ws.Range("J21").Value = 1.043694755208449
$ws.Range("K21").Value = 1.047080728779325
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.054015212425212
$ws.Range("N21").Value = 1.018447406067662
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036458887417296
$ws.Range("D22").Value = 1.043016305475683
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.049846133628521
$ws.Range("I22").Value = 1.039499903549702
$ws.Range("J22").Value = 1.043285081212946
$ws.Range("K22").Value = 1.046698510166574
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.05350252226006
$ws.Range("N22").Value = 1.018309634595711
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036807727029498
$ws.Range("D23").Value = 1.04328773381532
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.050186410756263
$ws.Range("I23").Value = 1.039587676733437
$ws.Range("J23").Value = 1.043502287161751
$ws.Range("K23").Value = 1.046901165745482
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.053774293969589
$ws.Range("N23").Value = 1.018382685215763
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038181821756701
$ws.Range("D24").Value = 1.044357082085894
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.051527806892744
$ws.Range("I24").Value = 1.039931001186444
$ws.Range("J24").Value = 1.044356736192872
$ws.Range("K24").Value = 1.047698236871176
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.05484456114872
$ws.Range("N24").Value = 1.018669936159764
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039777893503703
$ws.Range("D25").Value = 1.045599527555775
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.053088011713633
$ws.Range("I25").Value = 1.040324772069046
$ws.Range("J25").Value = 1.045346855236557
$ws.Range("K25").Value = 1.04862157486413
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.056087182796263
$ws.Range("N25").Value = 1.019267250563567
